# Scheduled market-data refresh: updates currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) across all item-sourcing sheets with freshly polled marketboard data.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (75 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1702
$ws.Range("I9").Value = 2494
$ws.Range("J9").Value = 712
$ws.Range("K9").Value = 2494
$ws.Range("L9").Value = 712
$ws.Range("M9").Value = -2325
$ws.Range("N9").Value = -1050
$ws.Range("H11").Value = 66.42856999999999
$ws.Range("I11").Value = 66.42856999999999
$ws.Range("K11").Value = 66.42856999999999
$ws.Range("M11").Value = 73.57143000000001
$ws.Range("H17").Value = 994.5
$ws.Range("J17").Value = 994.5
$ws.Range("L17").Value = 2983.5
$ws.Range("N17").Value = -3319.5
$ws.Range("H19").Value = 2471.6924
$ws.Range("J19").Value = 2113.6667
$ws.Range("L19").Value = 2113.6667
$ws.Range("N19").Value = -2463.6667
$ws.Range("H40").Value = 1796.9584
$ws.Range("I40").Value = 1478.9333
$ws.Range("K40").Value = 1478.9333
$ws.Range("M40").Value = -1303.9333
$ws.Range("I62").Value = 1575.2858
$ws.Range("J62").Value = 2405.75
$ws.Range("K62").Value = 1575.2858
$ws.Range("L62").Value = 2405.75
$ws.Range("M62").Value = -951.2858000000001
$ws.Range("N62").Value = -3653.75
$ws.Range("I65").Value = 1575.2858
$ws.Range("J65").Value = 2405.75
$ws.Range("K65").Value = 7876.429
$ws.Range("L65").Value = 12028.75
$ws.Range("M65").Value = -4756.429
$ws.Range("N65").Value = -18268.75
$ws.Range("H70").Value = 1850
$ws.Range("J70").Value = 1850
$ws.Range("L70").Value = 5550
$ws.Range("N70").Value = -6090
$ws.Range("H73").Value = 1850
$ws.Range("J73").Value = 1850
$ws.Range("L73").Value = 5550
$ws.Range("N73").Value = -7422
$ws.Range("H112").Value = 1545
$ws.Range("J112").Value = 1490
$ws.Range("L112").Value = 4470
$ws.Range("N112").Value = -6686
$ws.Range("H113").Value = 3503.625
$ws.Range("I113").Value = 3200
$ws.Range("J113").Value = 3604.8333
$ws.Range("K113").Value = 3200
$ws.Range("L113").Value = 3604.8333
$ws.Range("M113").Value = 54
$ws.Range("N113").Value = -10112.8333
$ws.Range("H132").Value = 1370.1875
$ws.Range("I132").Value = 1301.9166
$ws.Range("J132").Value = 1575
$ws.Range("K132").Value = 3905.7498
$ws.Range("L132").Value = 4725
$ws.Range("M132").Value = -1375.7498
$ws.Range("N132").Value = -9785
$ws.Range("H137").Value = 1494
$ws.Range("I137").Value = 1120.3636
$ws.Range("J137").Value = 1836.5
$ws.Range("K137").Value = 3361.0908
$ws.Range("L137").Value = 5509.5
$ws.Range("M137").Value = -811.0907999999999
$ws.Range("N137").Value = -10609.5
$ws.Range("H141").Value = 2041.3478
$ws.Range("I141").Value = 1916.0555
$ws.Range("J141").Value = 2492.4
$ws.Range("K141").Value = 5748.166499999999
$ws.Range("L141").Value = 7477.200000000001
$ws.Range("M141").Value = -568.1664999999994
$ws.Range("N141").Value = -17837.2

# ---- Sheet: ARM (32 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 192052.25
$ws.Range("I6").Value = 192052.25
$ws.Range("K6").Value = 192052.25
$ws.Range("M6").Value = -191879.25
$ws.Range("H32").Value = 2336124.5
$ws.Range("I32").Value = 3892805.8
$ws.Range("K32").Value = 3892805.8
$ws.Range("M32").Value = -3892518.8
$ws.Range("H45").Value = 1083.5714
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H61").Value = 2264.5715
$ws.Range("I61").Value = 2264.5715
$ws.Range("K61").Value = 2264.5715
$ws.Range("M61").Value = -2052.5715
$ws.Range("H74").Value = 5524.5
$ws.Range("I74").Value = 5524.5
$ws.Range("K74").Value = 5524.5
$ws.Range("M74").Value = -4650.5
$ws.Range("H77").Value = 5524.5
$ws.Range("I77").Value = 5524.5
$ws.Range("K77").Value = 27622.5
$ws.Range("M77").Value = -23254.5
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H136").Value = 2264.5715
$ws.Range("I136").Value = 2264.5715
$ws.Range("K136").Value = 6793.7145
$ws.Range("M136").Value = -4243.7145

# ---- Sheet: BSM (32 cell updates) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 115597.836
$ws.Range("J76").Value = 115597.836
$ws.Range("L76").Value = 115597.836
$ws.Range("N76").Value = -116227.836
$ws.Range("H79").Value = 115597.836
$ws.Range("J79").Value = 115597.836
$ws.Range("L79").Value = 115597.836
$ws.Range("N79").Value = -117781.836
$ws.Range("H80").Value = 688.4286
$ws.Range("J80").Value = 500.5
$ws.Range("L80").Value = 500.5
$ws.Range("N80").Value = -2496.5
$ws.Range("H83").Value = 688.4286
$ws.Range("J83").Value = 500.5
$ws.Range("L83").Value = 2502.5
$ws.Range("N83").Value = -12486.5
$ws.Range("H105").Value = 2133.8
$ws.Range("I105").Value = 2115.923
$ws.Range("K105").Value = 2115.923
$ws.Range("M105").Value = -368.9229999999998
$ws.Range("H107").Value = 2962.75
$ws.Range("J107").Value = 1997.5
$ws.Range("L107").Value = 1997.5
$ws.Range("N107").Value = -5837.5
$ws.Range("H137").Value = 78000
$ws.Range("J137").Value = 78000
$ws.Range("L137").Value = 78000
$ws.Range("N137").Value = -88200
$ws.Range("H138").Value = 56664.332
$ws.Range("J138").Value = 56664.332
$ws.Range("L138").Value = 56664.332
$ws.Range("N138").Value = -66944.33199999999

# ---- Sheet: CRP (31 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 17
$ws.Range("I10").Value = 17
$ws.Range("K10").Value = 17
$ws.Range("M10").Value = 122
$ws.Range("H31").Value = 2522.75
$ws.Range("I31").Value = 2380.3333
$ws.Range("K31").Value = 2380.3333
$ws.Range("M31").Value = -2085.3333
$ws.Range("H34").Value = 2522.75
$ws.Range("I34").Value = 2380.3333
$ws.Range("K34").Value = 2380.3333
$ws.Range("M34").Value = -2178.3333
$ws.Range("H105").Value = 2628.4119
$ws.Range("I105").Value = 1544.75
$ws.Range("J105").Value = 3591.6667
$ws.Range("K105").Value = 1544.75
$ws.Range("L105").Value = 3591.6667
$ws.Range("M105").Value = 202.25
$ws.Range("N105").Value = -7085.6667
$ws.Range("H107").Value = 531.6667
$ws.Range("I107").Value = 500
$ws.Range("K107").Value = 500
$ws.Range("M107").Value = 1420
$ws.Range("H117").Value = 48000
$ws.Range("J117").Value = 48000
$ws.Range("L117").Value = 48000
$ws.Range("N117").Value = -57178
$ws.Range("H134").Value = 3050.4
$ws.Range("I134").Value = 2950.8572
$ws.Range("K134").Value = 8852.571599999999
$ws.Range("M134").Value = -6317.571599999999

# ---- Sheet: CUL (19 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 12223069
$ws.Range("I4").Value = 12223069
$ws.Range("K4").Value = 36669207
$ws.Range("M4").Value = -36669095
$ws.Range("H97").Value = 1082.1428
$ws.Range("J97").Value = 362.33334
$ws.Range("L97").Value = 1087.00002
$ws.Range("N97").Value = -2079.00002
$ws.Range("H113").Value = 1659.6666
$ws.Range("I113").Value = 1247.25
$ws.Range("J113").Value = 1989.6
$ws.Range("K113").Value = 3741.75
$ws.Range("L113").Value = 5968.799999999999
$ws.Range("M113").Value = -1571.75
$ws.Range("N113").Value = -10308.8
$ws.Range("H137").Value = 2348.8333
$ws.Range("J137").Value = 4249.5
$ws.Range("L137").Value = 12748.5
$ws.Range("N137").Value = -22948.5

# ---- Sheet: GSM (4 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2831.3333
$ws.Range("I132").Value = 2831.3333
$ws.Range("K132").Value = 8493.999899999999
$ws.Range("M132").Value = -5963.999899999999

# ---- Sheet: LTW (11 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3233.1667
$ws.Range("I40").Value = 452
$ws.Range("J40").Value = 4623.75
$ws.Range("K40").Value = 452
$ws.Range("L40").Value = 4623.75
$ws.Range("M40").Value = -316
$ws.Range("N40").Value = -4895.75
$ws.Range("H93").Value = 899.6
$ws.Range("I93").Value = 874.75
$ws.Range("K93").Value = 874.75
$ws.Range("M93").Value = 373.25

# ---- Sheet: WVR (29 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 8549.5
$ws.Range("I14").Value = 18650
$ws.Range("J14").Value = 3499.25
$ws.Range("K14").Value = 18650
$ws.Range("L14").Value = 3499.25
$ws.Range("M14").Value = -18482
$ws.Range("N14").Value = -3835.25
$ws.Range("H81").Value = 1430770.4
$ws.Range("I81").Value = 1347.5
$ws.Range("K81").Value = 2695
$ws.Range("M81").Value = -1634
$ws.Range("H84").Value = 1430770.4
$ws.Range("I84").Value = 1347.5
$ws.Range("K84").Value = 13475
$ws.Range("M84").Value = -8171
$ws.Range("H96").Value = 1615
$ws.Range("I96").Value = 1700
$ws.Range("J96").Value = 1572.5
$ws.Range("K96").Value = 1700
$ws.Range("L96").Value = 1572.5
$ws.Range("M96").Value = -327
$ws.Range("N96").Value = -4318.5
$ws.Range("H132").Value = 2942.5557
$ws.Range("I132").Value = 3701.4
$ws.Range("J132").Value = 1994
$ws.Range("K132").Value = 11104.2
$ws.Range("L132").Value = 5982
$ws.Range("M132").Value = -8574.200000000001
$ws.Range("N132").Value = -11042

Write-Output "Updated 227 cells, added 4, cleared 2 across 8 sheets."
